# Regenerate merged AHB files
#
# 1) Rename the "_old"/"_new" comparison-column headers (row 1) to the
#    concrete format-version tags used by the new merge run:
#      *_old -> *_FV2410
#      *_new -> *_FV2504
# 2) Turn the used range A1:U72 into an Excel table ("Table1") so the
#    renamed headers become the table's column names.
# 3) Freeze the header row (row 1) so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) rename the header cells -------------------------------------------------
# Cells.Replace matches whole cell contents here, so longer/shorter names
# (e.g. "Segment_old" vs "Segmentname_old" / "Segment ID_old") don't collide.
$ws.Cells.Replace("Segmentname_old", "Segmentname_FV2410")
$ws.Cells.Replace("Segmentgruppe_old", "Segmentgruppe_FV2410")
$ws.Cells.Replace("Segment_old", "Segment_FV2410")
$ws.Cells.Replace("Datenelement_old", "Datenelement_FV2410")
$ws.Cells.Replace("Segment ID_old", "Segment ID_FV2410")
$ws.Cells.Replace("Code_old", "Code_FV2410")
$ws.Cells.Replace("Qualifier_old", "Qualifier_FV2410")
$ws.Cells.Replace("Beschreibung_old", "Beschreibung_FV2410")
$ws.Cells.Replace("Bedingungsausdruck_old", "Bedingungsausdruck_FV2410")
$ws.Cells.Replace("Bedingung_old", "Bedingung_FV2410")

$ws.Cells.Replace("Segmentname_new", "Segmentname_FV2504")
$ws.Cells.Replace("Segmentgruppe_new", "Segmentgruppe_FV2504")
$ws.Cells.Replace("Segment_new", "Segment_FV2504")
$ws.Cells.Replace("Datenelement_new", "Datenelement_FV2504")
$ws.Cells.Replace("Segment ID_new", "Segment ID_FV2504")
$ws.Cells.Replace("Code_new", "Code_FV2504")
$ws.Cells.Replace("Qualifier_new", "Qualifier_FV2504")
$ws.Cells.Replace("Beschreibung_new", "Beschreibung_FV2504")
$ws.Cells.Replace("Bedingungsausdruck_new", "Bedingungsausdruck_FV2504")
$ws.Cells.Replace("Bedingung_new", "Bedingung_FV2504")

# --- 2) convert the used range into a table --------------------------------------
$lo = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $ws.Range("A1:U72"),
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$lo.Name = "Table1"

# --- 3) freeze the header row ----------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
